$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.298.89"
Set-TextValue $ws.Range("E2") "  +1.27%  "
Set-TextValue $ws.Range("D3") "1.831.48"
Set-TextValue $ws.Range("E3") "  +0.71%  "
Set-TextValue $ws.Range("E4") "  +0.75%  "
Set-TextValue $ws.Range("D5") "314.06"
Set-TextValue $ws.Range("E5") "  +1.25%  "
Set-TextValue $ws.Range("D6") "1.009"
Set-TextValue $ws.Range("E6") "  +0.65%  "
Set-TextValue $ws.Range("D7") "0.4737"
Set-TextValue $ws.Range("E7") "  +1.78%  "
Set-TextValue $ws.Range("E8") "  +0.49%  "
Set-TextValue $ws.Range("D9") "0.07445"
Set-TextValue $ws.Range("E9") "  +1.21%  "
Set-TextValue $ws.Range("E10") "  +1.28%  "
Set-TextValue $ws.Range("E11") "  +1.11%  "
Set-TextValue $ws.Range("D12") "1.906.69"
Set-TextValue $ws.Range("E12") "  +3.73%  "
Set-TextValue $ws.Range("D13") "0.07304"
Set-TextValue $ws.Range("E13") "  +2.67%  "
Set-TextValue $ws.Range("D14") "5.419"
Set-TextValue $ws.Range("E14") "  +0.07%  "
Set-TextValue $ws.Range("D15") "93.78"
Set-TextValue $ws.Range("E15") "  +2.47%  "
Set-TextValue $ws.Range("D16") "6.551"
Set-TextValue $ws.Range("E16") "  +0.55%  "
Set-TextValue $ws.Range("E17") "  +0.31%  "
Set-TextValue $ws.Range("D18") "0.000008780"
Set-TextValue $ws.Range("E18") "  +0.85%  "
Set-TextValue $ws.Range("D20") "27.669.30"
Set-TextValue $ws.Range("E20") "  +2.58%  "
Set-TextValue $ws.Range("D21") "14.75"
Set-TextValue $ws.Range("E21") "  +0.60%  "
Set-TextValue $ws.Range("D22") "5.281"
Set-TextValue $ws.Range("E22") "  -0.26%  "
Set-TextValue $ws.Range("E23") "  +0.61%  "
Set-TextValue $ws.Range("D24") "2.106.98"
Set-TextValue $ws.Range("E24") "  +3.44%  "
Set-TextValue $ws.Range("D25") "1.897"
Set-TextValue $ws.Range("E25") "  +0.37%  "
Set-TextValue $ws.Range("D26") "151.69"
Set-TextValue $ws.Range("E26") "  +0.69%  "
Set-TextValue $ws.Range("E27") "  +0.66%  "
Set-TextValue $ws.Range("D28") "2.131"
Set-TextValue $ws.Range("E28") "  -0.42%  "
Set-TextValue $ws.Range("D29") "5.225"
Set-TextValue $ws.Range("E29") "  -0.51%  "
Set-TextValue $ws.Range("D30") "117.15"
Set-TextValue $ws.Range("E30") "  +0.48%  "
Set-TextValue $ws.Range("D31") "0.08985"
Set-TextValue $ws.Range("E31") "  +1.02%  "
Set-TextValue $ws.Range("B32") "ARBITRUM"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D32") "1.174"
Set-TextValue $ws.Range("E32") "  +0.90%  "
Set-TextValue $ws.Range("B33") "ImmutableX"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D33") "0.7482"
Set-TextValue $ws.Range("E33") "  -1.51%  "
Set-TextValue $ws.Range("E34") "  +0.72%  "
Set-TextValue $ws.Range("D35") "2.949"
Set-TextValue $ws.Range("E35") "  +1.38%  "
Set-TextValue $ws.Range("E36") "  +0.78%  "
Set-TextValue $ws.Range("D37") "1.094"
Set-TextValue $ws.Range("E37") "  -0.06%  "
Set-TextValue $ws.Range("D38") "0.05336"
Set-TextValue $ws.Range("E38") "  +0.78%  "
Set-TextValue $ws.Range("E39") "  +0.43%  "
Set-TextValue $ws.Range("D40") "2.429"
Set-TextValue $ws.Range("E40") "  +3.70%  "
Set-TextValue $ws.Range("D41") "2.957"
Set-TextValue $ws.Range("E41") "  -0.44%  "
Set-TextValue $ws.Range("D42") "7.226"
Set-TextValue $ws.Range("E42") "  +0.59%  "
Set-TextValue $ws.Range("D43") "0.5285"
Set-TextValue $ws.Range("E43") "  -0.15%  "
Set-TextValue $ws.Range("E44") "  +0.21%  "
Set-TextValue $ws.Range("D45") "8.486"
Set-TextValue $ws.Range("E45") "  +0.47%  "
Set-TextValue $ws.Range("D46") "0.4908"
Set-TextValue $ws.Range("D47") "10.51"
Set-TextValue $ws.Range("E47") "  +0.25%  "
Set-TextValue $ws.Range("D48") "105.01"
Set-TextValue $ws.Range("E48") "  +1.61%  "
Set-TextValue $ws.Range("E49") "  +0.73%  "
Set-TextValue $ws.Range("D50") "1.663"
Set-TextValue $ws.Range("E50") "  -0.36%  "
Set-TextValue $ws.Range("D51") "0.06296"
Set-TextValue $ws.Range("E51") "  +0.08%  "
